# Insert a new weekly price-report row for "Vega Monumental Concepción" /
# Mandarina at row 54, pushing the existing rows 54-108 down to 55-109
# (dimension grows from A1:T108 to A1:T109).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 54..108 down one row, leaving a blank row 54 to populate.
$ws.Rows.Item(54).Insert()

$ws.Range("A54").Value = 11
$ws.Range("B54").Value = "Vega Monumental Concepción"
$ws.Range("C54").Value = "Bíobío"
$ws.Range("D54").Value = 44721
$ws.Range("E54").Value = 8
$ws.Range("F54").Value = "Fruta"
$ws.Range("G54").Value = 100102
$ws.Range("H54").Value = "Cítricos"
$ws.Range("I54").Value = 100102004
$ws.Range("J54").Value = "Mandarina"
$ws.Range("K54").Value = "Clemenuless"
$ws.Range("L54").Value = "Primera"
$ws.Range("M54").Value = 300
$ws.Range("N54").Value = 7000
$ws.Range("O54").Value = 8000
$ws.Range("P54").Value = 7500
$ws.Range("Q54").Value = "$/caja 15 kilos"
$ws.Range("R54").Value = "Provincia de Limarí"
$ws.Range("S54").Value = 500
$ws.Range("T54").Value = 15
